$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New (reversed) order of period values for E16:E60
$periods = @(
    "2003","2002","2001",
    "1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901",
    "1812","1811","1810","1809","1808","1807","1806","1805","1804","1803","1802","1801",
    "1712","1711","1710","1709","1708","1707","1706","1705","1704","1703","1702","1701",
    "1612","1611","1610","1609","1608","1607"
)

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
}
